# Update horarios workbook: Linea 141 schedule refresh (scrape run 04:17:03)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:17:03"
$ws1.Range("A3").Value = "Total filas: 25"

$rows1 = @(
    @(16, "04:17:03", "04:31", "215_ALUAR",          14,  "LP1912"),
    @(17, "02:49:45", "04:35", "215_ALUAR",         106,  "LP1912"),
    @(18, "03:00:18", "04:44", "215_ALUAR",         104,  "LP1912"),
    @(19, "03:42:43", "04:45", "215A_EL PATO",        63, "LP1912"),
    @(20, "04:17:03", "04:53", "11_ETCHEVERRY",        36, "LP1912"),
    @(21, "04:17:03", "05:16", "17_ROMERO",            59, "LP1912"),
    @(22, "04:17:03", "05:22", "23_HERNANDEZ",          65, "LP1912"),
    @(23, "03:42:43", "05:34", "215B_EL PATO",        112, "LP1912"),
    @(24, "03:42:43", "05:35", "14_ABASTO",            113, "LP1912"),
    @(25, "04:17:03", "05:35", "215B_EL PATO",          78, "LP1912"),
    @(26, "04:17:03", "05:36", "14_ABASTO",             79, "LP1912"),
    @(27, "04:17:03", "05:46", "15_ABASTO",             89, "LP1912"),
    @(28, "04:17:03", "06:05", "16_SANTA ANA",         108, "LP1912"),
    @(29, "04:17:03", "06:12", "215A_EL PATO",         115, "LP1912"),
    @(30, "04:17:03", "06:14", "225_HARAS DEL SUR",    117, "LP1912")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:17:03"
$ws2.Range("A3").Value = "Total filas: 12"

$rows2 = @(
    @(11, "04:17:03", "04:31", "215_ALUAR",        14,  "LP1912"),
    @(12, "02:49:45", "04:35", "215_ALUAR",        106, "LP1912"),
    @(13, "03:00:18", "04:44", "215_ALUAR",        104, "LP1912"),
    @(14, "03:42:43", "04:45", "215A_EL PATO",       63, "LP1912"),
    @(15, "03:42:43", "05:34", "215B_EL PATO",      112, "LP1912"),
    @(16, "04:17:03", "05:35", "215B_EL PATO",       78, "LP1912"),
    @(17, "04:17:03", "06:12", "215A_EL PATO",      115, "LP1912")
)

foreach ($row in $rows2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:17:03"
$ws3.Range("A3").Value = "Total filas: 2"

$ws3.Cells.Item(7, 1).Value = "04:17:03"
$ws3.Cells.Item(7, 2).Value = "05:44"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 87
$ws3.Cells.Item(7, 5).Value = "L6173"
